$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "... не ясно, кто играет [ первую важную роль][_GoBack][? ]"
#           -> "... не ясно, кто играет первую важную роль? " (bookmark removed)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " первую важную роль? ", $true, $false, $false, $false, $false,
    $true, 1, $false, " первую важную роль? ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "». В частности, имеется ли доступ к градиенту в " + "«" -> merged
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "». В частности, имеется ли доступ к градиенту в «", $true, $false, $false, $false, $false,
    $true, 1, $false, "». В частности, имеется ли доступ к градиенту в «", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: "»" + "? Возможно, в виде псевдокода ... методов." -> merged
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "»? Возможно, в виде псевдокода стоит вставить и алгоритмы для этих «обычных» методов.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "»? Возможно, в виде псевдокода стоит вставить и алгоритмы для этих «обычных» методов.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: split the closing paragraph into two paragraphs:
#   1) "К коду замечаний нет."
#   2) original remark, trimmed, with the _GoBack bookmark re-homed right
#      after "читателю"
# ---------------------------------------------------------------------------

# Find the paragraph that holds the long closing remark.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Отмечу, что все указанные недочеты*результатами работы*") {
        $target = $p
    }
}

# Insert a new (empty) paragraph right after it - same technique Word uses
# when you press Enter at the end of a paragraph.
$target.Range.InsertParagraphAfter() | Out-Null

# Re-resolve the paragraph index, since the paragraph collection changed.
$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Отмечу, что все указанные недочеты*результатами работы*") {
        $targetIdx = $idx
    }
}
$newPara = $d.Paragraphs.Item($targetIdx + 1)

# Populate the new paragraph with the original remark, minus the trailing
# " разобраться с результатами работы" (kept ending: "...читателю.").
$newPara.Range.InsertBefore(
    "Отмечу, что все указанные недочеты не несут критического характера. Их основной посыл – сообщить о местах в работе, которые могут быть непонятны не погруженному в тематику читателю.") | Out-Null

# Replace the original paragraph's text with the short remark.
$target.Range.Text = "К коду замечаний нет."

# Re-home the _GoBack bookmark: drop it right after "читателю" (before the
# final period) in the newly created paragraph.
$r = $d.Content
$r.Find.Execute("тематику читателю.", $true, $false, $false, $false, $false,
                $true, 1, $false, "", 0) | Out-Null
$bmPos = $r.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Host "Edit complete"
